# QUBES Code First commit
# New test-data rows were appended to the generator's shared string pool and
# the two data cells on the "TestData" sheet (C5, C6) were advanced to the
# next pair of generated values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")

$ws.Range("C5").Value = "Var1-VS1P320220126"
$ws.Range("C6").Value = "WPL031076"
